# Update Querymodel ARCtrl reference:
#  - rename the "MS" assay worksheet to "MS_Cold"
#  - record a stray annotation value ("s") in column I below the table,
#    matching the re-saved workbook produced by the newer Excel build

$wb = $excel.ActiveWorkbook

$msSheet = $wb.Worksheets.Item("MS")
$msSheet.Name = "MS_Cold"

# Make it the active sheet/selection, as it was when the value was entered
$msSheet.Activate()

$msSheet.Range("I45").Value = "s"
$msSheet.Range("I45").Select() | Out-Null
